$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("J14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("J14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1

$ws.Range("H15").Value = 0

$ws.Range("J14").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 1

$ws.Range("K15").Value = 0

$ws.Range("L15").Value = 0

$ws.Range("M15").Value = 0

$ws.Range("N15").Value = -66.666666666666

# Row 16
$ws.Range("C16").Value = 4

$ws.Range("D16").Value = 3

$ws.Range("E16").Value = 33.333333333333

$ws.Range("F16").Value = 6

$ws.Range("G16").Value = 6

$ws.Range("H16").Value = 0

$ws.Range("I16").Value = 21

$ws.Range("J16").Value = 21

$ws.Range("K16").Value = 0

$ws.Range("L16").Value = 40

$ws.Range("M16").Value = 16.666666666666

$ws.Range("N16").Value = -76.666666666666

# Row 17
$ws.Range("C17").Value = 4

$ws.Range("D17").Value = 3

$ws.Range("E17").Value = 33.333333333333

$ws.Range("F17").Value = 10

$ws.Range("G17").Value = 10

$ws.Range("H17").Value = 0

$ws.Range("I17").Value = 26

$ws.Range("J17").Value = 30

$ws.Range("K17").Value = -13.333333333333

$ws.Range("L17").Value = 44.444444444444

$ws.Range("M17").Value = 36.842105263157

$ws.Range("N17").Value = -61.764705882352

# Row 18
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("F18").Value = 7

$ws.Range("H18").Value = 75

$ws.Range("I18").Value = 22

$ws.Range("K18").Value = 4.761904761904

$ws.Range("L18").Value = -33.333333333333

$ws.Range("M18").Value = -21.428571428571

$ws.Range("N18").Value = -79.245283018867

# Row 19
$ws.Range("C19").Value = 3

$ws.Range("J14").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = 200

$ws.Range("F19").Value = 15

$ws.Range("G19").Value = 6

$ws.Range("H19").Value = 150

$ws.Range("I19").Value = 45

$ws.Range("J19").Value = 40

$ws.Range("K19").Value = 12.5

$ws.Range("L19").Value = 2.272727272727

$ws.Range("M19").Value = -6.25

$ws.Range("N19").Value = -18.181818181818

# Row 20
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4163)

$ws.Range("D20").Value = 2

$ws.Range("E20").Value = -100

$ws.Range("F20").Value = 6

$ws.Range("G20").Value = 8

$ws.Range("H20").Value = -25

$ws.Range("I20").Value = 16

$ws.Range("J20").Value = 19

$ws.Range("K20").Value = -15.78947368421

$ws.Range("L20").Value = 14.285714285714

$ws.Range("M20").Value = 45.454545454545

$ws.Range("N20").Value = -86.324786324786

# Row 21
$ws.Range("C21").Value = 14

$ws.Range("E21").Value = 55.555555555555

$ws.Range("F21").Value = 45

$ws.Range("G21").Value = 35

$ws.Range("H21").Value = 28.571428571428

$ws.Range("I21").Value = 131

$ws.Range("J21").Value = 133

$ws.Range("K21").Value = -1.503759398496

$ws.Range("L21").Value = 3.968253968253

$ws.Range("M21").Value = 3.968253968253

$ws.Range("N21").Value = -70.294784580498

# Row 23
$ws.Range("C23").Value = 4

$ws.Range("D23").Value = 2

$ws.Range("E23").Value = 100

$ws.Range("F23").Value = 12

$ws.Range("G23").Value = 9

$ws.Range("H23").Value = 33.333333333333

$ws.Range("I23").Value = 29

$ws.Range("J23").Value = 35

$ws.Range("K23").Value = -17.142857142857

$ws.Range("L23").Value = 0

$ws.Range("M23").Value = 107.142857142857

# Row 24
$ws.Range("C24").Value = 4

$ws.Range("D24").Value = 17

$ws.Range("E24").Value = -76.470588235294

$ws.Range("F24").Value = 36

$ws.Range("G24").Value = 64

$ws.Range("H24").Value = -43.75

$ws.Range("I24").Value = 118

$ws.Range("J24").Value = 179

$ws.Range("K24").Value = -34.078212290502

$ws.Range("L24").Value = 31.111111111111

$ws.Range("M24").Value = 6.306306306306

# Row 25
$ws.Range("C14").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C25").PasteSpecial(-4163)

$ws.Range("E25").Value = -100

$ws.Range("F25").Value = 13

$ws.Range("H25").Value = -66.666666666666

$ws.Range("J25").Value = 102

$ws.Range("K25").Value = -61.764705882352

$ws.Range("L25").Value = 56

# Row 26
$ws.Range("C26").Value = 8

$ws.Range("D26").Value = 1

$ws.Range("E26").Value = 700

$ws.Range("F26").Value = 14

$ws.Range("G26").Value = 14

$ws.Range("H26").Value = 0

$ws.Range("I26").Value = 34

$ws.Range("J26").Value = 50

$ws.Range("K26").Value = -32

$ws.Range("L26").Value = -27.659574468085

$ws.Range("M26").Value = -49.253731343283

# Row 27
$ws.Range("J14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("I27").Value = 4

$ws.Range("K27").Value = 100

$ws.Range("L27").Value = 100

# Row 29
$ws.Range("J14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("J14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = 0

$ws.Range("J14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1

$ws.Range("J14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1

$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 0

$ws.Range("I29").Value = 2

$ws.Range("J29").Value = 2

$ws.Range("L29").Value = 100

$ws.Range("M29").Value = 0

$ws.Range("N29").Value = -77.777777777777

# Row 30
$ws.Range("J14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

$ws.Range("J14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = 0

$ws.Range("J14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

$ws.Range("J14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$ws.Range("K14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = 0

$ws.Range("I30").Value = 2

$ws.Range("J30").Value = 2

$ws.Range("L30").Value = 100

$ws.Range("M30").Value = 0

$ws.Range("N30").Value = -75

$excel.CutCopyMode = 0
Write-Host "cell edits done"